$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column D header
$ws.Range("D1").Value = "Correct_answer"

# Rows 2-51 correspond to the "Purple" colour group -> answer "l"
$ws.Range("D2:D51").Value = "l"

# Rows 52-101 correspond to the "Blue" colour group -> answer "s"
$ws.Range("D52:D101").Value = "s"

# Match the author's final viewport/selection state
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 83
$win.ScrollColumn = 1
$ws.Range("D52:D101").Select()
